# Adds a "Type of house" segmentation column (E) to the variable
# description sheet, classifying each variable's importance as
# High / Med / Low (one row uses the lowercase variant "high").
# Also stamps the header for the new column in G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label living in G2 (mirrors the column header row pattern already
# used elsewhere on row 2 for this sheet).
$ws.Range("G2").Value = "Type of house"

# Row -> comment value, row 2 is MSSubClass ... row 80 is SaleCondition.
$comments = @{
    2  = "High"; 3  = "High"; 4  = "Med";  5  = "High"; 6  = "Med";
    7  = "Low";  8  = "Low";  9  = "Med";  10 = "High"; 11 = "Med";
    12 = "Med";  13 = "High"; 14 = "Med";  15 = "Med";  16 = "Med";
    17 = "High"; 18 = "High"; 19 = "high"; 20 = "High"; 21 = "Med";
    22 = "Low";  23 = "Low";  24 = "Low";  25 = "Low";  26 = "Med";
    27 = "Med";  28 = "High"; 29 = "Med";  30 = "Low";  31 = "Low";
    32 = "Low";  33 = "Low";  34 = "Low";  35 = "Low";  36 = "Low";
    37 = "Low";  38 = "Low";  39 = "Med";  40 = "Med";  41 = "Med";
    42 = "Low";  43 = "Low";  44 = "Med";  45 = "Med";  46 = "Low";
    47 = "Low";  48 = "Low";  49 = "Low";  50 = "Med";  51 = "Low";
    52 = "Med";  53 = "Med";  54 = "Med";  55 = "Med";  56 = "Low";
    57 = "Low";  58 = "Low";  59 = "Low";  60 = "Low";  61 = "Low";
    62 = "High"; 63 = "High"; 64 = "Low";  65 = "Med";  66 = "Low";
    67 = "Low";  68 = "Low";  69 = "Low";  70 = "Low";  71 = "Low";
    72 = "Low";  73 = "Med";  74 = "Low";  75 = "Low";  76 = "Low";
    77 = "Med";  78 = "Med";  79 = "Med";  80 = "Med";
}

foreach ($row in ($comments.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 5).Value = $comments[$row]
}

# Leave the view scrolled/selected near the bottom of the data (mirrors
# where the editor's cursor ended up after filling the column down to
# row 80).
$ws.Range("E80").Select()
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
